$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.218.26'
$ws.Range("E2").Value = '  -1.09%  '
$ws.Range("D3").Value = '1.556.93'
$ws.Range("E3").Value = '  -0.90%  '
$ws.Range("E4").Value = '  -0.25%  '
$ws.Range("E5").Value = '  -0.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '288.29'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.05%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3808'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.91%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3316'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.60%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '44.68'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -7.47%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.146'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.17%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07403'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.20%  '
$ws.Range("E12").Value = '  -0.28%  '
$ws.Range("E13").Value = '  -3.40%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.843'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.25%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.750'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.45%  '
$ws.Range("D16").Value = '1.548.96'
$ws.Range("E16").Value = '  -1.58%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001075'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.78%  '
$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06654'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.68%  '
$ws.Range("B19").Value = 'Litecoin'
$ws.Range("C19").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '86.45'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.16%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.400'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.11%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.000'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '16.15'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.90%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.73'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.06%  '
$ws.Range("D24").Value = '22.236.61'
$ws.Range("E24").Value = '  -0.95%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.270'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.91%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.559'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.52%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '151.16'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.92%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.26'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.36%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.943'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.41%  '
$ws.Range("E30").Value = '  -0.81%  '
$ws.Range("D31").Value = '1.726.57'
$ws.Range("E31").Value = '  -1.32%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.091'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.63%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.914'
$ws.Range("D33").Style = "Normal"
$ws.Range("E34").Value = '  -5.35%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.331'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.70%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08212'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.27%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06327'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.92%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02336'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.05%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.319'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.28%  '
$ws.Range("E40").Value = '  -4.92%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.234'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '10.99'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.82%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6069'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.77%  '
$ws.Range("E44").Value = '  -0.30%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.76'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.21%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.746'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.99%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5870'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.43%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '122.33'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.28%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.967'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.40%  '
$ws.Range("E50").Value = '  -3.05%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07053'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.95%  '
